$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.990.40"
$ws.Range("E2").Value = "  +1.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.273.57"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.46"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.23"
$ws.Range("E6").Value = "  +5.01%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("E9").Value = "  +3.93%  "

$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.846.50"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.63"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.040.28"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("E16").Value = "  +2.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.270.63"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.86"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "381.85"
$ws.Range("E20").Value = "  +2.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.72"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.45"
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  +0.76%  "

$ws.Range("E25").Value = "  +1.70%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.190"
$ws.Range("E26").Value = "  +6.00%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.82"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.86"
$ws.Range("E29").Value = "  +4.41%  "

$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  +6.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.87"
$ws.Range("E32").Value = "  +1.25%  "

$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.28"
$ws.Range("E34").Value = "  +1.15%  "

$ws.Range("E35").Value = "  +2.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.28"
$ws.Range("E36").Value = "  -3.02%  "

$ws.Range("E37").Value = "  -0.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.840"
$ws.Range("E38").Value = "  -1.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.82"
$ws.Range("E39").Value = "  +4.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.62"
$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("E41").Value = "  +4.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.61"
$ws.Range("E42").Value = "  +0.87%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0695"
$ws.Range("E43").Value = "  +2.99%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.43"
$ws.Range("E44").Value = "  +2.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.44"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.653.98"
$ws.Range("E46").Value = "  -3.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "342.74"
$ws.Range("E47").Value = "  -3.12%  "

$ws.Range("E48").Value = "  +1.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.02"
$ws.Range("E49").Value = "  +3.93%  "

$ws.Range("E50").Value = "  +1.65%  "

$ws.Range("E51").Value = "  -0.05%  "
